# Update "想去人数" (interest count) values in column F on both the
# "展览" sheet and the "全部类型" sheet, matching the refreshed data
# snapshot described in the commit ("Update gh-pages to output generated
# at 456a3b4").

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 309
$wsExhibit.Range("F7").Value = 941
$wsExhibit.Range("F8").Value = 345
$wsExhibit.Range("F11").Value = 1411
$wsExhibit.Range("F13").Value = 1317
$wsExhibit.Range("F14").Value = 2975
$wsExhibit.Range("F15").Value = 373
$wsExhibit.Range("F25").Value = 3415

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 309
$wsAll.Range("F17").Value = 941
$wsAll.Range("F18").Value = 345
$wsAll.Range("F21").Value = 1411
$wsAll.Range("F23").Value = 1317
$wsAll.Range("F24").Value = 2975
$wsAll.Range("F25").Value = 373
$wsAll.Range("F37").Value = 3415
